$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 0.1506834993896362
$ws.Cells.Item(2, 4).Value = 0.1005736726642574
$ws.Cells.Item(2, 5).Value = 0.1396915254104911
$ws.Cells.Item(2, 6).Value = 2.447612376704498
$ws.Cells.Item(2, 7).Value = 0.00249952884151275
$ws.Cells.Item(2, 9).Value = 1.746282104336686
$ws.Cells.Item(2, 10).Value = 0.2074473435006823
$ws.Cells.Item(2, 12).Value = 0.2117037805259443
$ws.Cells.Item(2, 15).Value = 6.773999207379518
# Row 3
$ws.Cells.Item(3, 3).Value = 0.1506001191999786
$ws.Cells.Item(3, 4).Value = 0.1007978339278388
$ws.Cells.Item(3, 5).Value = 0.1396102605362533
$ws.Cells.Item(3, 6).Value = 2.411058583748769
$ws.Cells.Item(3, 7).Value = 0.002504394110525603
$ws.Cells.Item(3, 9).Value = 1.719032899959544
$ws.Cells.Item(3, 10).Value = 0.2067835050900939
$ws.Cells.Item(3, 12).Value = 0.2115744371488333
$ws.Cells.Item(3, 15).Value = 6.656197380233436
# Row 4
$ws.Cells.Item(4, 3).Value = 0.1506033809576479
$ws.Cells.Item(4, 4).Value = 0.1009607465917526
$ws.Cells.Item(4, 5).Value = 0.1396130617121845
$ws.Cells.Item(4, 6).Value = 2.389927942158621
$ws.Cells.Item(4, 7).Value = 0.002507542966596385
$ws.Cells.Item(4, 9).Value = 1.703247198003993
$ws.Cells.Item(4, 10).Value = 0.2064653220329475
$ws.Cells.Item(4, 12).Value = 0.2115765548513266
$ws.Cells.Item(4, 15).Value = 6.587612858722423
# Row 5
$ws.Cells.Item(5, 3).Value = 0.1506184217761515
$ws.Cells.Item(5, 4).Value = 0.101033504592003
$ws.Cells.Item(5, 5).Value = 0.1396274676000129
$ws.Cells.Item(5, 6).Value = 2.381646935565854
$ws.Cells.Item(5, 7).Value = 0.002508866911946211
$ws.Cells.Item(5, 9).Value = 1.697051819144519
$ws.Cells.Item(5, 10).Value = 0.2063581490956423
$ws.Cells.Item(5, 12).Value = 0.2115979369358527
$ws.Cells.Item(5, 15).Value = 6.560604663502261
# Row 6
$ws.Cells.Item(6, 3).Value = 0.1506217479725152
$ws.Cells.Item(6, 4).Value = 0.1010459711268865
$ws.Cells.Item(6, 5).Value = 0.1396306612125624
$ws.Cells.Item(6, 6).Value = 2.380291794617975
$ws.Cells.Item(6, 7).Value = 0.002509089218230279
$ws.Cells.Item(6, 9).Value = 1.696037412083484
$ws.Cells.Item(6, 10).Value = 0.2063417115621533
$ws.Cells.Item(6, 12).Value = 0.2116027272131191
$ws.Cells.Item(6, 15).Value = 6.556176738048009
# Row 7
$ws.Cells.Item(7, 3).Value = 0.1506035282615343
$ws.Cells.Item(7, 4).Value = 0.1009617020225484
$ws.Cells.Item(7, 5).Value = 0.1396132022702954
$ws.Cells.Item(7, 6).Value = 2.389814926386819
$ws.Cells.Item(7, 7).Value = 0.002507560656430805
$ws.Cells.Item(7, 9).Value = 1.703162684060729
$ws.Cells.Item(7, 10).Value = 0.2064637855962204
$ws.Cells.Item(7, 12).Value = 0.2115767601141698
$ws.Cells.Item(7, 15).Value = 6.587244810478467
# Row 8
$ws.Cells.Item(8, 3).Value = 0.1506434575370008
$ws.Cells.Item(8, 4).Value = 0.1006457254054602
$ws.Cells.Item(8, 5).Value = 0.1396525751937077
$ws.Cells.Item(8, 6).Value = 2.434735744743548
$ws.Cells.Item(8, 7).Value = 0.002501172935840412
$ws.Cells.Item(8, 9).Value = 1.736690154748217
$ws.Cells.Item(8, 10).Value = 0.207199895637622
$ws.Cells.Item(8, 12).Value = 0.2116422695768563
$ws.Cells.Item(8, 15).Value = 6.732602645548866
# Row 9
$ws.Cells.Item(9, 3).Value = 0.151153226531477
$ws.Cells.Item(9, 4).Value = 0.1002260900631669
$ws.Cells.Item(9, 5).Value = 0.1401475027033285
$ws.Cells.Item(9, 6).Value = 2.533273707637704
$ws.Cells.Item(9, 7).Value = 0.002489922381638247
$ws.Cells.Item(9, 9).Value = 1.809958868860903
$ws.Cells.Item(9, 10).Value = 0.2093530118710731
$ws.Cells.Item(9, 12).Value = 0.212417194507573
$ws.Cells.Item(9, 15).Value = 7.047465294632445
# Row 10
$ws.Cells.Item(10, 3).Value = 0.1517901743448462
$ws.Cells.Item(10, 4).Value = 0.1000390000918685
$ws.Cells.Item(10, 5).Value = 0.1407654421211042
$ws.Cells.Item(10, 6).Value = 2.612086996391838
$ws.Cells.Item(10, 7).Value = 0.002482425708685732
$ws.Cells.Item(10, 9).Value = 1.86841092062312
$ws.Cells.Item(10, 10).Value = 0.2113681746740355
$ws.Cells.Item(10, 12).Value = 0.2133803393407021
$ws.Cells.Item(10, 15).Value = 7.297133943381652
# Row 11
$ws.Cells.Item(11, 3).Value = 0.1521368116122517
$ws.Cells.Item(11, 4).Value = 0.09998006562654993
$ws.Cells.Item(11, 5).Value = 0.1411017313729133
$ws.Cells.Item(11, 6).Value = 2.649345923571218
$ws.Cells.Item(11, 7).Value = 0.002479180446187943
$ws.Cells.Item(11, 9).Value = 1.896014281248512
$ws.Cells.Item(11, 10).Value = 0.2123791936029633
$ws.Cells.Item(11, 12).Value = 0.2139039825385893
$ws.Cells.Item(11, 15).Value = 7.414734569018833
# Row 12
$ws.Cells.Item(12, 3).Value = 0.1522762418226904
$ws.Cells.Item(12, 4).Value = 0.09996150008998228
$ws.Cells.Item(12, 5).Value = 0.1412370034829316
$ws.Cells.Item(12, 6).Value = 2.663657843630517
$ws.Cells.Item(12, 7).Value = 0.002477975138173772
$ws.Cells.Item(12, 9).Value = 1.906613189990537
$ws.Cells.Item(12, 10).Value = 0.2127756097622182
$ws.Cells.Item(12, 12).Value = 0.2141145604211445
$ws.Cells.Item(12, 15).Value = 7.459848082973167
# Row 13
$ws.Cells.Item(13, 3).Value = 0.1522458500575041
$ws.Cells.Item(13, 4).Value = 0.09996533184454748
$ws.Cells.Item(13, 5).Value = 0.1412075178116226
$ws.Cells.Item(13, 6).Value = 2.660566488450911
$ws.Cells.Item(13, 7).Value = 0.002478233675135444
$ws.Cells.Item(13, 9).Value = 1.904324020137494
$ws.Cells.Item(13, 10).Value = 0.2126896310700275
$ws.Cells.Item(13, 12).Value = 0.2140686624520427
$ws.Cells.Item(13, 15).Value = 7.450106222564784
# Row 14
$ws.Cells.Item(14, 3).Value = 0.1521481190347558
$ws.Cells.Item(14, 4).Value = 0.09997846310876923
$ws.Cells.Item(14, 5).Value = 0.1411127014784697
$ws.Cells.Item(14, 6).Value = 2.650519306838873
$ws.Cells.Item(14, 7).Value = 0.00247908081258978
$ws.Cells.Item(14, 9).Value = 1.896883329077042
$ws.Cells.Item(14, 10).Value = 0.2124115351092897
$ws.Cells.Item(14, 12).Value = 0.2139210607408728
$ws.Cells.Item(14, 15).Value = 7.418434435628569
# Row 15
$ws.Cells.Item(15, 3).Value = 0.1520893190325765
$ws.Cells.Item(15, 4).Value = 0.09998699460826188
$ws.Cells.Item(15, 5).Value = 0.1410556557639389
$ws.Cells.Item(15, 6).Value = 2.644391541257392
$ws.Cells.Item(15, 7).Value = 0.002479602778810827
$ws.Cells.Item(15, 9).Value = 1.892344731923473
$ws.Cells.Item(15, 10).Value = 0.2122429600104638
$ws.Cells.Item(15, 12).Value = 0.213832250094427
$ws.Cells.Item(15, 15).Value = 7.399110234349394
# Row 16
$ws.Cells.Item(16, 3).Value = 0.1517686642811782
$ws.Cells.Item(16, 4).Value = 0.1000433772432068
$ws.Cells.Item(16, 5).Value = 0.1407445745287106
$ws.Cells.Item(16, 6).Value = 2.609680383338571
$ws.Cells.Item(16, 7).Value = 0.002482641103935677
$ws.Cells.Item(16, 9).Value = 1.866627400008113
$ws.Cells.Item(16, 10).Value = 0.2113040002657769
$ws.Cells.Item(16, 12).Value = 0.2133478377320728
$ws.Cells.Item(16, 15).Value = 7.289529623125077
# Row 17
$ws.Cells.Item(17, 3).Value = 0.151586512897822
$ws.Cells.Item(17, 4).Value = 0.1000846612311044
$ws.Cells.Item(17, 5).Value = 0.1405678641227261
$ws.Cells.Item(17, 6).Value = 2.588746795704651
$ws.Cells.Item(17, 7).Value = 0.002484547191223748
$ws.Cells.Item(17, 9).Value = 1.851110453494613
$ws.Cells.Item(17, 10).Value = 0.210752136425441
$ws.Cells.Item(17, 12).Value = 0.2130725588642903
$ws.Cells.Item(17, 15).Value = 7.22333773475367
# Row 18
$ws.Cells.Item(18, 3).Value = 0.1514870994834325
$ws.Cells.Item(18, 4).Value = 0.1001108709940368
$ws.Cells.Item(18, 5).Value = 0.1404714198056318
$ws.Cells.Item(18, 6).Value = 2.576838690617961
$ws.Cells.Item(18, 7).Value = 0.002485659061528936
$ws.Cells.Item(18, 9).Value = 1.842280860830556
$ws.Cells.Item(18, 10).Value = 0.210443595716086
$ws.Cells.Item(18, 12).Value = 0.2129222743105643
$ws.Cells.Item(18, 15).Value = 7.185644656867225
# Row 19
$ws.Cells.Item(19, 3).Value = 0.1514543600112432
$ws.Cells.Item(19, 4).Value = 0.1001201687878819
$ws.Cells.Item(19, 5).Value = 0.1404396579128999
$ws.Cells.Item(19, 6).Value = 2.572829524297418
$ws.Cells.Item(19, 7).Value = 0.002486038194331421
$ws.Cells.Item(19, 9).Value = 1.839307675527934
$ws.Cells.Item(19, 10).Value = 0.2103406536367345
$ws.Cells.Item(19, 12).Value = 0.2128727731993862
$ws.Cells.Item(19, 15).Value = 7.17294742090013
# Row 20
$ws.Cells.Item(20, 3).Value = 0.1516053490788636
$ws.Cells.Item(20, 4).Value = 0.1000800115211717
$ws.Cells.Item(20, 5).Value = 0.1405861376677997
$ws.Cells.Item(20, 6).Value = 2.590961511249617
$ws.Cells.Item(20, 7).Value = 0.002484342677818013
$ws.Cells.Item(20, 9).Value = 1.852752388202674
$ws.Cells.Item(20, 10).Value = 0.2108099645391661
$ws.Cells.Item(20, 12).Value = 0.213101029846456
$ws.Cells.Item(20, 15).Value = 7.230344762001778
# Row 21
$ws.Cells.Item(21, 3).Value = 0.1521766034872343
$ws.Cells.Item(21, 4).Value = 0.09997450441429834
$ws.Cells.Item(21, 5).Value = 0.1411403362758534
$ws.Cells.Item(21, 6).Value = 2.653464900058452
$ws.Cells.Item(21, 7).Value = 0.002478831348221717
$ws.Cells.Item(21, 9).Value = 1.899064871532005
$ws.Cells.Item(21, 10).Value = 0.212492850446786
$ws.Cells.Item(21, 12).Value = 0.2139640815894595
$ws.Cells.Item(21, 15).Value = 7.42772143151376
# Row 22
$ws.Cells.Item(22, 3).Value = 0.1525975438850793
$ws.Cells.Item(22, 4).Value = 0.09992741150873385
$ws.Cells.Item(22, 5).Value = 0.1415487359092253
$ws.Cells.Item(22, 6).Value = 2.695496654110912
$ws.Cells.Item(22, 7).Value = 0.002475366889443447
$ws.Cells.Item(22, 9).Value = 1.930184627031267
$ws.Cells.Item(22, 10).Value = 0.213671784344946
$ws.Cells.Item(22, 12).Value = 0.2145997413838074
$ws.Cells.Item(22, 15).Value = 7.560104227502904
# Row 23
$ws.Cells.Item(23, 3).Value = 0.1523685298818265
$ws.Cells.Item(23, 4).Value = 0.09995054944637971
$ws.Cells.Item(23, 5).Value = 0.1413265410027051
$ws.Cells.Item(23, 6).Value = 2.672955164908501
$ws.Cells.Item(23, 7).Value = 0.00247720339543811
$ws.Cells.Item(23, 9).Value = 1.913497342145476
$ws.Cells.Item(23, 10).Value = 0.2130353288472975
$ws.Cells.Item(23, 12).Value = 0.2142539289422984
$ws.Cells.Item(23, 15).Value = 7.4891386027316
# Row 24
$ws.Cells.Item(24, 3).Value = 0.1515968167072757
$ws.Cells.Item(24, 4).Value = 0.1000821059450061
$ws.Cells.Item(24, 5).Value = 0.1405778601597696
$ws.Cells.Item(24, 6).Value = 2.589959843127303
$ws.Cells.Item(24, 7).Value = 0.002484435088262795
$ws.Cells.Item(24, 9).Value = 1.85200978513798
$ws.Cells.Item(24, 10).Value = 0.2107837932467902
$ws.Cells.Item(24, 12).Value = 0.2130881332684211
$ws.Cells.Item(24, 15).Value = 7.227175759008787
# Row 25
$ws.Cells.Item(25, 3).Value = 0.1509691936218687
$ws.Cells.Item(25, 4).Value = 0.1003182753267673
$ws.Cells.Item(25, 5).Value = 0.1399689281676082
$ws.Cells.Item(25, 6).Value = 2.505493110964196
$ws.Cells.Item(25, 7).Value = 0.002492830271187865
$ws.Cells.Item(25, 9).Value = 1.789328805967472
$ws.Cells.Item(25, 10).Value = 0.2086945016658532
$ws.Cells.Item(25, 12).Value = 0.2121383810426991
$ws.Cells.Item(25, 15).Value = 6.959077938267114
